$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-string header text (Volume number, report date range) ---
$hdr1 = $ws.Range("A8")
$hdr1.Characters(21, 2).Text = "49"

$hdr2 = $ws.Range("C9")
$hdr2.Characters(27, 10).Text = "12/5/2022"
$hdr2.Characters(47, 9).Text = "12/11/2022"

# --- Update crime statistics table (rows 15-29) ---
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial()
$excel.CutCopyMode = $false
$ws.Range("D15").Copy()
$ws.Range("G15").PasteSpecial()
$excel.CutCopyMode = $false
$ws.Range("E15").Copy()
$ws.Range("H15").PasteSpecial()
$excel.CutCopyMode = $false
$ws.Range("M15").Value = 33.333333333333
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -54.545454545454
$ws.Range("I16").Value = 102
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = 7.368421052631
$ws.Range("L16").Value = -12.068965517241
$ws.Range("M16").Value = -54.867256637168
$ws.Range("N16").Value = -90.981432360742
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 54.545454545454
$ws.Range("I17").Value = 165
$ws.Range("J17").Value = 132
$ws.Range("K17").Value = 25
$ws.Range("L17").Value = -24.311926605504
$ws.Range("M17").Value = 28.90625
$ws.Range("N17").Value = -64.362850971922
$ws.Range("D15").Copy()
$ws.Range("C18").PasteSpecial()
$excel.CutCopyMode = $false
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -53.333333333333
$ws.Range("I18").Value = 112
$ws.Range("J18").Value = 97
$ws.Range("K18").Value = 15.463917525773
$ws.Range("L18").Value = -13.178294573643
$ws.Range("M18").Value = -36
$ws.Range("N18").Value = -89.282296650717
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 16.666666666666
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 6.896551724137
$ws.Range("I19").Value = 323
$ws.Range("J19").Value = 277
$ws.Range("K19").Value = 16.606498194945
$ws.Range("L19").Value = 9.121621621621
$ws.Range("M19").Value = -28.381374722838
$ws.Range("N19").Value = -38.358778625954
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 120
$ws.Range("I20").Value = 92
$ws.Range("J20").Value = 58
$ws.Range("K20").Value = 58.620689655172
$ws.Range("L20").Value = 37.31343283582
$ws.Range("M20").Value = -5.154639175257
$ws.Range("N20").Value = -87.042253521126
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 71
$ws.Range("H21").Value = 1.408450704225
$ws.Range("I21").Value = 806
$ws.Range("J21").Value = 669
$ws.Range("K21").Value = 20.478325859491
$ws.Range("L21").Value = -4.389086595492
$ws.Range("M21").Value = -26.190476190476
$ws.Range("N21").Value = -79.470198675496
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 17
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = 70
$ws.Range("L22").Value = -34.615384615384
$ws.Range("M22").Value = -67.307692307692
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = -20
$ws.Range("I23").Value = 83
$ws.Range("J23").Value = 103
$ws.Range("K23").Value = -19.417475728155
$ws.Range("L23").Value = -29.661016949152
$ws.Range("M23").Value = -18.627450980392
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 42.857142857142
$ws.Range("F24").Value = 50
$ws.Range("H24").Value = -20.63492063492
$ws.Range("I24").Value = 719
$ws.Range("J24").Value = 547
$ws.Range("K24").Value = 31.44424131627
$ws.Range("L24").Value = 8.939393939393
$ws.Range("M24").Value = -35.86083853702
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = -8
$ws.Range("I25").Value = 355
$ws.Range("J25").Value = 246
$ws.Range("K25").Value = 44.30894308943
$ws.Range("L25").Value = 44.30894308943
$ws.Range("M25").Value = 28.623188405797
$ws.Range("D15").Copy()
$ws.Range("C26").PasteSpecial()
$excel.CutCopyMode = $false
$ws.Range("D15").Copy()
$ws.Range("G26").PasteSpecial()
$excel.CutCopyMode = $false
$ws.Range("E15").Copy()
$ws.Range("H26").PasteSpecial()
$excel.CutCopyMode = $false
$ws.Range("D27").Value = 3
$ws.Range("D15").Copy()
$ws.Range("F27").PasteSpecial()
$excel.CutCopyMode = $false
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 43
$ws.Range("K27").Value = -37.209302325581
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("I28").Value = 9
$ws.Range("K28").Value = -18.181818181818
$ws.Range("L28").Value = -55
$ws.Range("M28").Value = -52.631578947368
$ws.Range("N28").Value = -92.173913043478
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("I29").Value = 9
$ws.Range("K29").Value = -10
$ws.Range("L29").Value = -40
$ws.Range("M29").Value = -47.058823529411
$ws.Range("N29").Value = -91.089108910891
